$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-11-01 Friday" "2024-11-02 Saturday"

Replace-Text "458÷3=152, 2" "386÷3=128, 2"
Replace-Text "141÷3=47, 0" "451÷8=56, 3"
Replace-Text "137÷2=68, 1" "660÷6=110, 0"
Replace-Text "329÷6=54, 5" "446÷3=148, 2"
Replace-Text "764÷2=382, 0" "133÷6=22, 1"

Replace-Text "750÷2=375, 0" "455÷4=113, 3"
Replace-Text "611÷5=122, 1" "343÷2=171, 1"
Replace-Text "845÷6=140, 5" "392÷4=98, 0"
Replace-Text "722÷5=144, 2" "695÷9=77, 2"
Replace-Text "222÷2=111, 0" "810÷4=202, 2"

Replace-Text "995÷6=165, 5" "303÷7=43, 2"
Replace-Text "999÷5=199, 4" "612÷8=76, 4"
Replace-Text "370÷4=92, 2" "960÷9=106, 6"
Replace-Text "884÷9=98, 2" "743÷5=148, 3"
Replace-Text "236÷5=47, 1" "409÷5=81, 4"

Replace-Text "782÷7=111, 5" "749÷3=249, 2"
Replace-Text "524÷8=65, 4" "943÷8=117, 7"
Replace-Text "809÷5=161, 4" "371÷5=74, 1"
Replace-Text "276÷4=69, 0" "265÷7=37, 6"
Replace-Text "514÷2=257, 0" "519÷2=259, 1"

Replace-Text "135÷5=27, 0" "875÷5=175, 0"
Replace-Text "314÷5=62, 4" "100÷9=11, 1"
Replace-Text "451÷4=112, 3" "356÷5=71, 1"
Replace-Text "608÷7=86, 6" "477÷6=79, 3"
Replace-Text "319÷9=35, 4" "296÷5=59, 1"
